$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cxcl12"
$ws.Cells.Item(2, 3).Value = "Cxcr3"
$ws.Cells.Item(2, 4).Value = "M1"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 145.002022
$ws.Cells.Item(2, 8).Value = 290.004044
$ws.Cells.Item(2, 9).Value = 0.2376076076529671
$ws.Cells.Item(2, 10).Value = 0.1752271928451109
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.7166990000000001
$ws.Cells.Item(2, 14).Value = 2.150097
$ws.Cells.Item(2, 15).Value = 0.2276207788704612
$ws.Cells.Item(2, 16).Value = 0.2276207788704611
$ws.Cells.Item(2, 17).Value = 103.922804165378
$ws.Cells.Item(2, 18).Value = 623.5368249922681
$ws.Cells.Item(2, 19).Value = 0.05408442871951533
$ws.Cells.Item(2, 20).Value = 0.03988535011468863

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cxcl12"
$ws.Cells.Item(3, 3).Value = "Cxcr3"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 145.002022
$ws.Cells.Item(3, 8).Value = 290.004044
$ws.Cells.Item(3, 9).Value = 0.2376076076529671
$ws.Cells.Item(3, 10).Value = 0.1752271928451109
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.431954666666666
$ws.Cells.Item(3, 14).Value = 7.295864
$ws.Cells.Item(3, 15).Value = 0.7723792211295388
$ws.Cells.Item(3, 16).Value = 0.7723792211295388
$ws.Cells.Item(3, 17).Value = 352.6383440790026
$ws.Cells.Item(3, 18).Value = 2115.830064474016
$ws.Cells.Item(3, 19).Value = 0.1835231789334518
$ws.Cells.Item(3, 20).Value = 0.1353418427304222

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Cxcl12"
$ws.Cells.Item(4, 3).Value = "Cxcr3"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 295.9742430000001
$ws.Cells.Item(4, 8).Value = 887.9227290000001
$ws.Cells.Item(4, 9).Value = 0.4849982837213674
$ws.Cells.Item(4, 10).Value = 0.5365035780881736
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.7166990000000001
$ws.Cells.Item(4, 14).Value = 2.150097
$ws.Cells.Item(4, 15).Value = 0.2276207788704612
$ws.Cells.Item(4, 16).Value = 0.2276207788704611
$ws.Cells.Item(4, 17).Value = 212.1244439838571
$ws.Cells.Item(4, 18).Value = 1909.119995854713
$ws.Cells.Item(4, 19).Value = 0.1103956870914946
$ws.Cells.Item(4, 20).Value = 0.1221193623112194

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cxcl12"
$ws.Cells.Item(5, 3).Value = "Cxcr3"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 295.9742430000001
$ws.Cells.Item(5, 8).Value = 887.9227290000001
$ws.Cells.Item(5, 9).Value = 0.4849982837213674
$ws.Cells.Item(5, 10).Value = 0.5365035780881736
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.431954666666666
$ws.Cells.Item(5, 14).Value = 7.295864
$ws.Cells.Item(5, 15).Value = 0.7723792211295388
$ws.Cells.Item(5, 16).Value = 0.7723792211295388
$ws.Cells.Item(5, 17).Value = 719.7959414769841
$ws.Cells.Item(5, 18).Value = 6478.163473292857
$ws.Cells.Item(5, 19).Value = 0.3746025966298729
$ws.Cells.Item(5, 20).Value = 0.4143842157769543

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Cxcl12"
$ws.Cells.Item(6, 3).Value = "Cxcr3"
$ws.Cells.Item(6, 4).Value = "M1"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.1807913333333333
$ws.Cells.Item(6, 8).Value = 0.542374
$ws.Cells.Item(6, 9).Value = 0.0002962537736040913
$ws.Cells.Item(6, 10).Value = 0.0003277149938370313
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.7166990000000001
$ws.Cells.Item(6, 14).Value = 2.150097
$ws.Cells.Item(6, 15).Value = 0.2276207788704612
$ws.Cells.Item(6, 16).Value = 0.2276207788704611
$ws.Cells.Item(6, 17).Value = 0.1295729678086667
$ws.Cells.Item(6, 18).Value = 1.166156710278
$ws.Cells.Item(6, 19).Value = 0.00006743351469107653
$ws.Cells.Item(6, 20).Value = 0.00007459474214471345

# Row 7
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Cxcl12"
$ws.Cells.Item(7, 3).Value = "Cxcr3"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.1807913333333333
$ws.Cells.Item(7, 8).Value = 0.542374
$ws.Cells.Item(7, 9).Value = 0.0002962537736040913
$ws.Cells.Item(7, 10).Value = 0.0003277149938370313
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.431954666666666
$ws.Cells.Item(7, 14).Value = 7.295864
$ws.Cells.Item(7, 15).Value = 0.7723792211295388
$ws.Cells.Item(7, 16).Value = 0.7723792211295388
$ws.Cells.Item(7, 17).Value = 0.4396763267928889
$ws.Cells.Item(7, 18).Value = 3.957086941136
$ws.Cells.Item(7, 19).Value = 0.0002288202589130147
$ws.Cells.Item(7, 20).Value = 0.0002531202516923179

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Cxcl12"
$ws.Cells.Item(8, 3).Value = "Cxcr3"
$ws.Cells.Item(8, 4).Value = "M1"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.738772
$ws.Cells.Item(8, 8).Value = 2.216316
$ws.Cells.Item(8, 9).Value = 0.001210588963518025
$ws.Cells.Item(8, 10).Value = 0.00133914970902166
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.7166990000000001
$ws.Cells.Item(8, 14).Value = 2.150097
$ws.Cells.Item(8, 15).Value = 0.2276207788704612
$ws.Cells.Item(8, 16).Value = 0.2276207788704611
$ws.Cells.Item(8, 17).Value = 0.5294771536280001
$ws.Cells.Item(8, 18).Value = 4.765294382652
$ws.Cells.Item(8, 19).Value = 0.0002755552027679571
$ws.Cells.Item(8, 20).Value = 0.0003048182997916617

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Cxcl12"
$ws.Cells.Item(9, 3).Value = "Cxcr3"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.738772
$ws.Cells.Item(9, 8).Value = 2.216316
$ws.Cells.Item(9, 9).Value = 0.001210588963518025
$ws.Cells.Item(9, 10).Value = 0.00133914970902166
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.431954666666666
$ws.Cells.Item(9, 14).Value = 7.295864
$ws.Cells.Item(9, 15).Value = 0.7723792211295388
$ws.Cells.Item(9, 16).Value = 0.7723792211295388
$ws.Cells.Item(9, 17).Value = 1.796660013002666
$ws.Cells.Item(9, 18).Value = 16.169940117024
$ws.Cells.Item(9, 19).Value = 0.0009350337607500677
$ws.Cells.Item(9, 20).Value = 0.001034331409229999

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Cxcl12"
$ws.Cells.Item(10, 3).Value = "Cxcr3"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 137.6068803333333
$ws.Cells.Item(10, 8).Value = 412.820641
$ws.Cells.Item(10, 9).Value = 0.225489556501436
$ws.Cells.Item(10, 10).Value = 0.2494358391462612
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.7166990000000001
$ws.Cells.Item(10, 14).Value = 2.150097
$ws.Cells.Item(10, 15).Value = 0.2276207788704612
$ws.Cells.Item(10, 16).Value = 0.2276207788704611
$ws.Cells.Item(10, 17).Value = 98.62271352801967
$ws.Cells.Item(10, 18).Value = 887.604421752177
$ws.Cells.Item(10, 19).Value = 0.05132610847801172
$ws.Cells.Item(10, 20).Value = 0.05677677998467905

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutro"
$ws.Cells.Item(11, 2).Value = "Cxcl12"
$ws.Cells.Item(11, 3).Value = "Cxcr3"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 137.6068803333333
$ws.Cells.Item(11, 8).Value = 412.820641
$ws.Cells.Item(11, 9).Value = 0.225489556501436
$ws.Cells.Item(11, 10).Value = 0.2494358391462612
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 2.431954666666666
$ws.Cells.Item(11, 14).Value = 7.295864
$ws.Cells.Item(11, 15).Value = 0.7723792211295388
$ws.Cells.Item(11, 16).Value = 0.7723792211295388
$ws.Cells.Item(11, 17).Value = 334.6536947920915
$ws.Cells.Item(11, 18).Value = 3011.883253128824
$ws.Cells.Item(11, 19).Value = 0.1741634480234242
$ws.Cells.Item(11, 20).Value = 0.1926590591615822

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Cxcl12"
$ws.Cells.Item(12, 3).Value = "Cxcr3"
$ws.Cells.Item(12, 4).Value = "M1"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 30.755622
$ws.Cells.Item(12, 8).Value = 61.511244
$ws.Cells.Item(12, 9).Value = 0.05039770938710747
$ws.Cells.Item(12, 10).Value = 0.03716652521759548
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.7166990000000001
$ws.Cells.Item(12, 14).Value = 2.150097
$ws.Cells.Item(12, 15).Value = 0.2276207788704612
$ws.Cells.Item(12, 16).Value = 0.2276207788704611
$ws.Cells.Item(12, 17).Value = 22.042523531778
$ws.Cells.Item(12, 18).Value = 132.255141190668
$ws.Cells.Item(12, 19).Value = 0.01147156586398056
$ws.Cells.Item(12, 20).Value = 0.008459873417937717

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Cxcl12"
$ws.Cells.Item(13, 3).Value = "Cxcr3"
$ws.Cells.Item(13, 4).Value = "M2"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 30.755622
$ws.Cells.Item(13, 8).Value = 61.511244
$ws.Cells.Item(13, 9).Value = 0.05039770938710747
$ws.Cells.Item(13, 10).Value = 0.03716652521759548
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.431954666666666
$ws.Cells.Item(13, 14).Value = 7.295864
$ws.Cells.Item(13, 15).Value = 0.7723792211295388
$ws.Cells.Item(13, 16).Value = 0.7723792211295388
$ws.Cells.Item(13, 17).Value = 74.796278449136
$ws.Cells.Item(13, 18).Value = 448.777670694816
$ws.Cells.Item(13, 19).Value = 0.03892614352312691
$ws.Cells.Item(13, 20).Value = 0.02870665179965776
